# Remove the "Today's Lesson" agenda slide from the deck.
$p = $ppt.ActivePresentation

$targetIndex = -1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.HasText) {
                if ($shape.TextFrame.TextRange.Text -eq "Today's Lesson") {
                    $targetIndex = $i
                }
            }
        }
    }
}

if ($targetIndex -eq -1) {
    # Fallback: "Today's Lesson" is the 2nd slide in the deck.
    $targetIndex = 2
}

$p.Slides.Item($targetIndex).Delete()
